$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.718.52'
$ws.Cells.Item(2, 5).Value = '  +3.08%  '
$ws.Cells.Item(3, 4).Value = '1.863.15'
$ws.Cells.Item(3, 5).Value = '  +2.97%  '
$ws.Cells.Item(4, 5).Value = '  +2.78%  '
$ws.Cells.Item(5, 4).Value = '''323.86'
$ws.Cells.Item(5, 5).Value = '  +3.83%  '
$ws.Cells.Item(6, 4).Value = '''1.032'
$ws.Cells.Item(6, 5).Value = '  +2.57%  '
$ws.Cells.Item(7, 4).Value = '''0.4411'
$ws.Cells.Item(7, 5).Value = '  +2.93%  '
$ws.Cells.Item(8, 4).Value = '''0.3798'
$ws.Cells.Item(8, 5).Value = '  +3.08%  '
$ws.Cells.Item(9, 4).Value = '''0.07447'
$ws.Cells.Item(9, 5).Value = '  +2.62%  '
$ws.Cells.Item(10, 4).Value = '''0.8832'
$ws.Cells.Item(10, 5).Value = '  +2.35%  '
$ws.Cells.Item(11, 4).Value = '''21.82'
$ws.Cells.Item(11, 5).Value = '  +3.13%  '
$ws.Cells.Item(12, 4).Value = '1.871.81'
$ws.Cells.Item(12, 5).Value = '  -8.64%  '
$ws.Cells.Item(13, 4).Value = '''5.553'
$ws.Cells.Item(13, 5).Value = '  +2.97%  '
$ws.Cells.Item(14, 4).Value = '''6.740'
$ws.Cells.Item(14, 5).Value = '  +1.84%  '
$ws.Cells.Item(15, 4).Value = '''0.07198'
$ws.Cells.Item(15, 5).Value = '  +3.91%  '
$ws.Cells.Item(16, 4).Value = '''83.97'
$ws.Cells.Item(16, 5).Value = '  +3.96%  '
$ws.Cells.Item(17, 4).Value = '''1.038'
$ws.Cells.Item(17, 5).Value = '  +2.46%  '
$ws.Cells.Item(18, 4).Value = '''0.000009090'
$ws.Cells.Item(18, 5).Value = '  +2.81%  '
$ws.Cells.Item(19, 4).Value = '''1.032'
$ws.Cells.Item(19, 5).Value = '  +2.57%  '
$ws.Cells.Item(20, 4).Value = '''15.53'
$ws.Cells.Item(20, 5).Value = '  +1.90%  '
$ws.Cells.Item(21, 4).Value = '27.742.93'
$ws.Cells.Item(21, 5).Value = '  +3.05%  '
$ws.Cells.Item(22, 4).Value = '''5.305'
$ws.Cells.Item(22, 5).Value = '  +2.24%  '
$ws.Cells.Item(23, 4).Value = '''11.43'
$ws.Cells.Item(23, 5).Value = '  +4.60%  '
$ws.Cells.Item(24, 4).Value = '''158.41'
$ws.Cells.Item(24, 5).Value = '  +2.81%  '
$ws.Cells.Item(25, 4).Value = '''1.937'
$ws.Cells.Item(25, 5).Value = '  +2.81%  '
$ws.Cells.Item(26, 4).Value = '''18.82'
$ws.Cells.Item(26, 5).Value = '  +2.52%  '
$ws.Cells.Item(27, 4).Value = '''1.989'
$ws.Cells.Item(27, 5).Value = '  +4.20%  '
$ws.Cells.Item(28, 4).Value = '''5.325'
$ws.Cells.Item(28, 5).Value = '  +1.85%  '
$ws.Cells.Item(29, 4).Value = '''117.57'
$ws.Cells.Item(29, 5).Value = '  +2.59%  '
$ws.Cells.Item(30, 4).Value = '''0.09086'
$ws.Cells.Item(30, 5).Value = '  +1.62%  '
$ws.Cells.Item(31, 4).Value = '''1.213'
$ws.Cells.Item(31, 5).Value = '  +4.84%  '
$ws.Cells.Item(32, 4).Value = '''0.7692'
$ws.Cells.Item(32, 5).Value = '  +4.05%  '
$ws.Cells.Item(33, 4).Value = '''3.007'
$ws.Cells.Item(33, 5).Value = '  +7.14%  '
$ws.Cells.Item(34, 4).Value = '''4.571'
$ws.Cells.Item(34, 5).Value = '  +3.33%  '
$ws.Cells.Item(35, 4).Value = '''1.033'
$ws.Cells.Item(35, 5).Value = '  +2.60%  '
$ws.Cells.Item(36, 4).Value = '''1.164'
$ws.Cells.Item(36, 5).Value = '  +3.80%  '
$ws.Cells.Item(37, 4).Value = '''0.01986'
$ws.Cells.Item(37, 5).Value = '  +3.36%  '
$ws.Cells.Item(38, 4).Value = '''0.05342'
$ws.Cells.Item(38, 5).Value = '  +2.39%  '
$ws.Cells.Item(39, 4).Value = '''0.5197'
$ws.Cells.Item(39, 5).Value = '  +2.37%  '
$ws.Cells.Item(40, 5).Value = '  +3.11%  '
$ws.Cells.Item(41, 4).Value = '''0.1691'
$ws.Cells.Item(41, 5).Value = '  +2.63%  '
$ws.Cells.Item(42, 4).Value = '''6.852'
$ws.Cells.Item(42, 5).Value = '  +6.28%  '
$ws.Cells.Item(43, 4).Value = '''8.710'
$ws.Cells.Item(43, 5).Value = '  +4.85%  '
$ws.Cells.Item(44, 4).Value = '''109.40'
$ws.Cells.Item(44, 5).Value = '  +2.06%  '
$ws.Cells.Item(45, 5).Value = '  +2.64%  '
$ws.Cells.Item(46, 4).Value = '''1.727'
$ws.Cells.Item(46, 5).Value = '  +5.05%  '
$ws.Cells.Item(47, 4).Value = '''0.4690'
$ws.Cells.Item(47, 5).Value = '  +2.47%  '
$ws.Cells.Item(48, 4).Value = '''0.06425'
$ws.Cells.Item(48, 5).Value = '  +2.45%  '
$ws.Cells.Item(49, 4).Value = '''1.873'
$ws.Cells.Item(49, 5).Value = '  +3.46%  '
$ws.Cells.Item(50, 4).Value = '''39.68'
$ws.Cells.Item(50, 5).Value = '  +4.59%  '
$ws.Cells.Item(51, 4).Value = '''64.41'
$ws.Cells.Item(51, 5).Value = '  +1.45%  '
